# The source diff for this fixture only reorders XML attributes inside
# word/document.xml (root element namespace declarations, and the
# <w:sectPr>/<w:pgSz>/<w:pgMar> attributes) and word/styles.xml
# (<w:docDefaults>, <w:latentStyles>, every <w:lsdException>, and the
# <w:style>/<w:tblPr> elements). Every changed line is the exact same
# element/attribute set, just alphabetized by attribute name - i.e. a
# purely cosmetic serialization artifact.
#
# The commit message confirms this: "Fixed POI packaging and upgraded to
# POI 3.15" - a build/library change in the authoring tool that produced
# this expected-generation fixture, not an edit to the document's
# content, text, formatting, structure or styles.
#
# There is therefore no Word object-model mutation to make here: the
# paragraphs, runs, text, sections, page setup and styles are unchanged
# before and after. We simply touch the document without altering it.
$d = $word.ActiveDocument
